# TxAuthorize.xlsx — "授權理由" (Reason) redesign on the DBD layout sheet.
#
# The "授權理由代碼/ReasonCode" field is dropped, "授權理由/Reason" is
# renamed to "交易理由/TradeReason", and "ReasonJson" is renamed to
# "ReasonFAJson". Deleting the ReasonCode row (entire-row delete, shifting
# everything below it up by one) reproduces the row renumbering, the
# dimension shrink (A1:G21 -> A1:G20) and the sheet1 "spans" clean-up for
# the trailing KeyID/CreateDate/... rows that Excel recomputes when it
# rewrites a sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 12 is "ReasonCode" / "授權理由代碼" — remove it outright, shifting
# rows 13:21 up into 12:20.
$ws.Range("A12:G12").EntireRow.Delete()

# The row that shifted up into 12 was "Reason" / "授權理由" -> rename it.
$ws.Range("C12").Value = "交易理由"
$ws.Range("B12").Value = "TradeReason"

# The row that shifted up into 13 was "ReasonJson" -> rename it.
$ws.Range("B13").Value = "ReasonFAJson"

# Last touched/selected cell in the saved file.
$ws.Range("B14").Select()
